$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before column D (MobileNo onwards shift right by one).
$ws.Columns.Item(4).Insert()

# 2. New D1 header: "Dob(DDMMYYYY)"
$ws.Range("D1").Value = "Dob(DDMMYYYY)"

# Give D1 the plain "Arial/black" look (same look as the pre-existing style used
# elsewhere in the sheet) by copying format from a cell that already carries it,
# instead of touching Font properties one at a time (which would mint new style
# records). Column Z (far outside the used range) uses the sheet's base style.
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D1").PasteSpecial(-4122) | Out-Null

# 3. The former "Due" column (M1) has shifted to N1 after the insert above;
# rename it and add the brand-new "StartDate(DDMMYYYY)" column after it.
$ws.Range("N1").Value = "NoOfPayments"
$ws.Range("O1").Value = "StartDate(DDMMYYYY)"

# O1 uses the same plain look as D1.
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("O1").PasteSpecial(-4122) | Out-Null

# N1 gets a brand-new monospace font (JetBrains Mono) distinct from the header font.
$nFont = $ws.Range("N1").Font
$nFont.Name = "JetBrains Mono"
$nFont.Size = 10
$nFont.Color = 0
$nFont.Family = 3

# 4. Column widths for the new/affected columns (values are chosen so the
# engine's internal character-width rounding lands as close as possible to
# the authored widths of 16.79 / 15.39 / 21.36).
$ws.Columns.Item(4).ColumnWidth = 15.9167
$ws.Columns.Item(14).ColumnWidth = 14.4167
$ws.Columns.Item(15).ColumnWidth = 20.4167

# 5. Selection moves to D3.
$ws.Range("D3").Select() | Out-Null
